$d = $word.ActiveDocument

# The "Recorded Classes Link" table is the first (only) table in the document.
$t = $d.Tables.Item(1)

# Add a brand-new row at the end of the table - Word clones the formatting
# (cell widths / run properties) of the row above automatically.
$newRow = $t.Rows.Add()

# --- Column 1: "17th June" with "th" superscripted ---------------------
$cell1 = $newRow.Cells.Item(1)
$r1 = $cell1.Range
$r1.Text = "17th June"
$r1Start = $cell1.Range.Start
$supRange = $d.Range($r1Start + 2, $r1Start + 4)
$supRange.Font.Superscript = $true

# --- Column 2: "DOM" -----------------------------------------------------
$cell2 = $newRow.Cells.Item(2)
$cell2.Range.Text = "DOM"

# --- Column 3: hyperlink to the recorded class video ---------------------
$cell3 = $newRow.Cells.Item(3)
$url = "https://www.youtube.com/watch?v=L4yMKwcc8aM"
$cell3.Range.Text = $url
$linkRange = $cell3.Range
$moveResult = $linkRange.MoveEnd(1, -1)
$d.Hyperlinks.Add($linkRange, $url, $null, $null, $url) | Out-Null

Write-Output "Added 17th June / DOM / $url row to the Recorded Classes table."
